$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("D5").Value = "신호 공간(signal space)"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2022/01/12/signal_space.html"

# Row 6
$ws.Range("D6").Value = "[Optimization] 최적화 알고리즘 :: GA(Genetic Algorithm, 유전 알고리즘)란? GA 예시, R로 GA 구현하기"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Optimization-%EC%B5%9C%EC%A0%81%ED%99%94-%EC%95%8C%EA%B3%A0%EB%A6%AC%EC%A6%98-GA%EB%9E%80"

# Row 26
$ws.Range("D26").Value = "2021 인공지능 경진대회 참가기"

# Row 37
$ws.Range("D37").Value = "[Paper Review] SituatedQA: Incorporating Extra-Linguistic Contexts into QA"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1904&mod=document&pageid=1"

# Row 51
$ws.Range("D51").Value = "[우분투] 서비스 실행 상태 확인 명령어"
$ws.Range("E51").Value = "https://bskyvision.com/1237"

# Row 52
$ws.Range("D52").Value = "3판 맛보기) R에서 정수형(integer)과 실수형(numeric)"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2377&utm_source=rss&utm_medium=rss&utm_campaign=r%25ec%2597%2590%25ec%2584%259c-%25ec%25a0%2595%25ec%2588%2598%25ed%2598%2595integer%25ea%25b3%25bc-%25ec%258b%25a4%25ec%2588%2598%25ed%2598%2595numeric"
